$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''68.487.73'
$ws.Range("E2").Value = '  +1.26%  '
$ws.Range("D3").Value = '''3.356.95'
$ws.Range("E3").Value = '  +0.99%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''584.87'
$ws.Range("E5").Value = '  +0.94%  '
$ws.Range("D6").Value = '''177.55'
$ws.Range("E6").Value = '  +1.53%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '''0.591'
$ws.Range("E8").Value = '  +0.55%  '
$ws.Range("D9").Value = '''0.185'
$ws.Range("E9").Value = '  +4.24%  '
$ws.Range("D10").Value = '''0.583'
$ws.Range("E10").Value = '  +1.42%  '
$ws.Range("D11").Value = '''48.02'
$ws.Range("E11").Value = '  +6.01%  '
$ws.Range("E12").Value = '  +2.22%  '
$ws.Range("D13").Value = '''693.10'
$ws.Range("E13").Value = '  +4.95%  '
$ws.Range("D14").Value = '''3.901.85'
$ws.Range("E14").Value = '  +0.87%  '
$ws.Range("D15").Value = '''8.48'
$ws.Range("E15").Value = '  +1.16%  '
$ws.Range("D16").Value = '''68.434.63'
$ws.Range("E16").Value = '  +1.16%  '
$ws.Range("E17").Value = '  +1.48%  '
$ws.Range("D18").Value = '''3.370.42'
$ws.Range("E18").Value = '  +1.36%  '
$ws.Range("D19").Value = '''17.52'
$ws.Range("E19").Value = '  +1.10%  '
$ws.Range("D20").Value = '''11.26'
$ws.Range("E20").Value = '  +3.02%  '
$ws.Range("D21").Value = '''0.897'
$ws.Range("E21").Value = '  +1.27%  '
$ws.Range("D22").Value = '''5.50'
$ws.Range("E22").Value = '  +2.59%  '
$ws.Range("D23").Value = '''16.97'
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("D24").Value = '''100.11'
$ws.Range("E24").Value = '  +1.31%  '
$ws.Range("E25").Value = '  +1.92%  '
$ws.Range("E26").Value = '  +2.01%  '
$ws.Range("D27").Value = '''9.54'
$ws.Range("E27").Value = '  +3.43%  '
$ws.Range("D28").Value = '''33.11'
$ws.Range("E28").Value = '  -0.78%  '
$ws.Range("D29").Value = '''8.56'
$ws.Range("E29").Value = '  +1.78%  '
$ws.Range("D30").Value = '''6.99'
$ws.Range("E30").Value = '  -3.60%  '
$ws.Range("D31").Value = '''11.11'
$ws.Range("E31").Value = '  +1.71%  '
$ws.Range("D32").Value = '''549.99'
$ws.Range("E32").Value = '  -3.50%  '
$ws.Range("E33").Value = '  +0.83%  '
$ws.Range("D34").Value = '''58.20'
$ws.Range("E34").Value = '  +2.92%  '
$ws.Range("B35").Value = 'Dai'
$ws.Range("C35").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D35").Value = '''0.998'
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("B36").Value = 'Maker'
$ws.Range("C36").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D36").Value = '''3.705.39'
$ws.Range("E36").Value = '  +0.99%  '
$ws.Range("E37").Value = '  +4.51%  '
$ws.Range("E38").Value = '  +8.90%  '
$ws.Range("D39").Value = '''34.78'
$ws.Range("E39").Value = '  +1.67%  '
$ws.Range("E40").Value = '  +2.77%  '
$ws.Range("E41").Value = '  +0.33%  '
$ws.Range("D42").Value = '''0.0₃0675'
$ws.Range("E42").Value = '  +2.29%  '
$ws.Range("D43").Value = '''0.336'
$ws.Range("E43").Value = '  +1.04%  '
$ws.Range("E44").Value = '  -3.18%  '
$ws.Range("E45").Value = '  +1.88%  '
$ws.Range("D46").Value = '''2.65'
$ws.Range("E46").Value = '  +2.25%  '
$ws.Range("D47").Value = '''0.129'
$ws.Range("E47").Value = '  +0.74%  '
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("E49").Value = '  -1.59%  '
$ws.Range("D50").Value = '''131.68'
$ws.Range("E50").Value = '  +1.66%  '
$ws.Range("D51").Value = '''2.63'
$ws.Range("E51").Value = '  -1.25%  '
